# Weekly update: insert a new data row at row 32 (pushing existing rows
# 32..59 down to 33..60) and populate it with the latest week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 32, shifting rows down.
$ws.Rows.Item(32).Insert(-4121)  # -4121 = xlShiftDown

# Populate the new row 32 with the new record's values.
$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C32").Value = "Los Lagos"
$ws.Range("D32").Value = 44957
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 100112030
$ws.Range("G32").Value = "Poroto granado"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 43000
$ws.Range("L32").Value = 43000
$ws.Range("M32").Value = 43000
$ws.Range("N32").Value = "$/saco 25 kilos"
$ws.Range("O32").Value = "Región del Maule"
$ws.Range("P32").Value = 1720
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
